# Applies the cryptos-list refresh described by the commit
# "Updated cryptos list on Mon Aug 19 22:08:52 UTC 2024 with GitHub Actions".
# Updates Coin (B), Link (C), Price (D) and Volume/1h (E) cells for the
# rows whose scraped values changed; two rows (Fetch.AI/SuiNetwork,
# Stellar/Hedera, Maker/RenderToken) also swapped rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (matches the source workbook,
# which stores every Coin/Link/Price/Volume cell as a string) without
# leaving a stray number-format style behind on the cell.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 "58.862.95"
Set-TextCell 2 5 "  -1.02%  "

# Row 3
Set-TextCell 3 4 "2.604.93"
Set-TextCell 3 5 "  -1.29%  "

# Row 4
Set-TextCell 4 5 "  +0.03%  "

# Row 5
Set-TextCell 5 4 "556.08"
Set-TextCell 5 5 "  +3.59%  "

# Row 6
Set-TextCell 6 4 "143.79"
Set-TextCell 6 5 "  -1.03%  "

# Row 7
Set-TextCell 7 5 "  -0.05%  "

# Row 8
Set-TextCell 8 4 "0.597"
Set-TextCell 8 5 "  +4.33%  "

# Row 9
Set-TextCell 9 4 "6.80"
Set-TextCell 9 5 "  -2.80%  "

# Row 10
Set-TextCell 10 5 "  -0.86%  "

# Row 11
Set-TextCell 11 4 "0.143"
Set-TextCell 11 5 "  +5.58%  "

# Row 12
Set-TextCell 12 4 "0.335"
Set-TextCell 12 5 "  -0.92%  "

# Row 13
Set-TextCell 13 4 "3.066.86"
Set-TextCell 13 5 "  -1.19%  "

# Row 14
Set-TextCell 14 4 "58.831.70"
Set-TextCell 14 5 "  -0.93%  "

# Row 15
Set-TextCell 15 4 "20.93"
Set-TextCell 15 5 "  -1.85%  "

# Row 16
Set-TextCell 16 4 "2.604.15"
Set-TextCell 16 5 "  -2.06%  "

# Row 17
Set-TextCell 17 4 "0.0000132"
Set-TextCell 17 5 "  -1.78%  "

# Row 18
Set-TextCell 18 5 "  -0.77%  "

# Row 19
Set-TextCell 19 4 "337.75"
Set-TextCell 19 5 "  +0.02%  "

# Row 20
Set-TextCell 20 4 "10.09"
Set-TextCell 20 5 "  -2.06%  "

# Row 21
Set-TextCell 21 4 "6.19"
Set-TextCell 21 5 "  -0.71%  "

# Row 22
Set-TextCell 22 4 "0.998"
Set-TextCell 22 5 "  -0.19%  "

# Row 23
Set-TextCell 23 4 "66.54"
Set-TextCell 23 5 "  +0.35%  "

# Row 24
Set-TextCell 24 4 "0.429"
Set-TextCell 24 5 "  +2.78%  "

# Row 25
Set-TextCell 25 5 "  -0.20%  "

# Row 26
Set-TextCell 26 5 "  -2.19%  "

# Row 27
Set-TextCell 27 4 "7.17"
Set-TextCell 27 5 "  -1.49%  "

# Row 28
Set-TextCell 28 4 "0.0₃0759"
Set-TextCell 28 5 "  +1.18%  "

# Row 29
Set-TextCell 29 5 "  -0.04%  "

# Row 30
Set-TextCell 30 5 "  +1.37%  "

# Row 31
Set-TextCell 31 4 "6.00"
Set-TextCell 31 5 "  +1.66%  "

# Row 32
Set-TextCell 32 4 "154.43"
Set-TextCell 32 5 "  +2.20%  "

# Row 33
Set-TextCell 33 4 "18.96"
Set-TextCell 33 5 "  +0.70%  "

# Row 34
Set-TextCell 34 5 "  -0.82%  "

# Row 35
Set-TextCell 35 2 "Fetch.AI"
Set-TextCell 35 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell 35 4 "0.923"
Set-TextCell 35 5 "  +10.23%  "

# Row 36
Set-TextCell 36 2 "SuiNetwork"
Set-TextCell 36 3 "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell 36 4 "0.906"
Set-TextCell 36 5 "  +7.60%  "

# Row 37
Set-TextCell 37 5 "  -0.33%  "

# Row 38
Set-TextCell 38 4 "37.13"
Set-TextCell 38 5 "  -0.56%  "

# Row 39
Set-TextCell 39 5 "  +0.41%  "

# Row 40
Set-TextCell 40 4 "3.60"
Set-TextCell 40 5 "  -0.37%  "

# Row 41
Set-TextCell 41 4 "283.28"
Set-TextCell 41 5 "  -0.64%  "

# Row 43
Set-TextCell 43 4 "0.599"
Set-TextCell 43 5 "  -0.32%  "

# Row 44
Set-TextCell 44 2 "Hedera"
Set-TextCell 44 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 44 4 "0.0537"
Set-TextCell 44 5 "  -0.10%  "

# Row 45
Set-TextCell 45 2 "Stellar"
Set-TextCell 45 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 45 4 "0.0953"
Set-TextCell 45 5 "  +1.16%  "

# Row 46
Set-TextCell 46 4 "10.61"
Set-TextCell 46 5 "  -1.27%  "

# Row 47
Set-TextCell 47 4 "0.0227"
Set-TextCell 47 5 "  +0.25%  "

# Row 48
Set-TextCell 48 2 "Maker"
Set-TextCell 48 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell 48 4 "1.950.79"
Set-TextCell 48 5 "  -0.52%  "

# Row 49
Set-TextCell 49 2 "RenderToken"
Set-TextCell 49 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 49 4 "4.56"
Set-TextCell 49 5 "  +0.02%  "

# Row 50
Set-TextCell 50 4 "118.66"
Set-TextCell 50 5 "  +6.50%  "

# Row 51
Set-TextCell 51 4 "17.94"
Set-TextCell 51 5 "  -2.47%  "
